$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "68.508.78"
$ws.Range("E2").Value = "  -1.69%  "

$ws.Range("D3").Value = "3.856.10"
$ws.Range("E3").Value = "  -1.00%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.58"
$ws.Range("E5").Value = "  -0.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.26"
$ws.Range("E6").Value = "  -1.33%  "

$ws.Range("D7").Value = "3.857.66"

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  -1.30%  "

$ws.Range("E10").Value = "  -1.95%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.49"
$ws.Range("E11").Value = "  +1.24%  "

$ws.Range("E12").Value = "  -2.06%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000267"
$ws.Range("E13").Value = "  +4.21%  "

$ws.Range("E14").Value = "  -3.03%  "

$ws.Range("D15").Value = "4.503.14"
$ws.Range("E15").Value = "  -1.10%  "

$ws.Range("D16").Value = "3.859.51"
$ws.Range("E16").Value = "  -1.00%  "

$ws.Range("D17").Value = "68.727.34"
$ws.Range("E17").Value = "  -1.37%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.56"
$ws.Range("E18").Value = "  -0.87%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.39"

$ws.Range("E20").Value = "  -0.82%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.22"
$ws.Range("E21").Value = "  +1.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "471.38"
$ws.Range("E22").Value = "  -3.98%  "

$ws.Range("E23").Value = "  -1.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000161"
$ws.Range("E24").Value = "  -2.87%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.52"
$ws.Range("E25").Value = "  -2.19%  "

$ws.Range("E26").Value = "  -2.88%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.15"
$ws.Range("E27").Value = "  -1.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.23"
$ws.Range("E28").Value = "  +0.80%  "

$ws.Range("E29").Value = "  +0.14%  "

$ws.Range("E30").Value = "  -0.68%  "

$ws.Range("D31").Value = "4.007.19"
$ws.Range("E31").Value = "  -1.03%  "

$ws.Range("E32").Value = "  -2.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.50"
$ws.Range("E33").Value = "  -1.39%  "

$ws.Range("E34").Value = "  -4.15%  "

$ws.Range("E35").Value = "  -3.05%  "

$ws.Range("D36").Value = "3.821.47"
$ws.Range("E36").Value = "  -1.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.104"
$ws.Range("E37").Value = "  -2.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.72"
$ws.Range("E38").Value = "  +8.26%  "

$ws.Range("E40").Value = "  -2.12%  "

$ws.Range("E42").Value = "  +0.10%  "

$ws.Range("E43").Value = "  -3.91%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.98"
$ws.Range("E44").Value = "  -5.21%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.75"
$ws.Range("E45").Value = "  +0.52%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000294"
$ws.Range("E46").Value = "  +6.99%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "417.09"
$ws.Range("E47").Value = "  -4.20%  "

$ws.Range("E49").Value = "  -1.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0360"
$ws.Range("E50").Value = "  -1.70%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "141.83"
$ws.Range("E51").Value = "  -0.43%  "

